$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "If you reach this point stop for a while!" paragraph, plus
#    the blank ListParagraph that immediately follows it.
# ---------------------------------------------------------------------------
$found = $d.Content.Find.Execute("If you reach this point stop for a while!")
if ($found) {
    $stopPara = $d.Content.Paragraphs.Item(1)
}

$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -eq "If you reach this point stop for a while!`r") {
        $nextPara = $paras.Item($i + 1)
        $delStart = $p.Range.Start
        $delEnd = $nextPara.Range.End
        $d.Range($delStart, $delEnd).Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 2) "Overwrite the doGet method." -> "Override the doGet() method."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Overwrite the doGet method.", $false, $false, $false, $false, $false, $true, 1, $false, "Override the doGet() method.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Both occurrences of "overwritten" -> "overridden"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("In the overwritten method build an html table with all the headers from the request.", $false, $false, $false, $false, $false, $true, 1, $false, "In the overridden method build an html table with all the headers from the request.", 2) | Out-Null

$d.Content.Find.Execute("In the overwritten method write to the response the following:", $false, $false, $false, $false, $false, $true, 1, $false, "In the overridden method write to the response the following:", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Move the "_GoBack" bookmark from the end of the "...and their values"
#    paragraph to right after "overridden" in the "build an html table..."
#    paragraph.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $old = $d.Bookmarks.Item("_GoBack")
    $old.Delete()
}

$rng = $d.Content
$rng.Find.Execute("In the overridden")
$newBmStart = $rng.End
$bmRange = $d.Range($newBmStart, $newBmStart)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
